$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("TestCase")
$ws2 = $wb.Worksheets.Item("Note")

# Insert a new column before column G (current TesterName column)
$ws1.Columns("G").Insert()

$ws1.Range("G1").Value = "Note"

$ws2.Range("B2").Value = "v1.3.1.0"
$ws2.Range("B3").Value = 43753
$ws2.Range("B4").Value = 18
$ws2.Range("C4").Value = "- Do not change ordinal of columns in the left of column System Validation. `n- Can change name of any columns.`n- Can add or remove columns in the right of column System Validation."
